# "Generate Report for Handback" - refresh the status/report rows for the
# 7ed2a854-01c6-4bd3-8381-0ee20e055f16.md file: it has now been handed back
# in sync with en-US (was "Ready for handoff"), with a fresh handback
# timestamp and no pending error.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 5 (7ed2a854...) zh-cn / de-de status columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E5").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F5").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 5 (7ed2a854...) Status / Latest Handback DateTime / Error Detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C5").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K5").Value = "2016-10-10 09:45:44"
$wsZhCn.Range("P5").Value = ""
$wsZhCn.Columns.Item(16).AutoFit()

# --- de-de sheet: row 5 (7ed2a854...) Status / Latest Handback DateTime / Error Detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C5").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K5").Value = "2016-10-10 09:46:00"
$wsDeDe.Range("P5").Value = ""
$wsDeDe.Columns.Item(16).AutoFit()

Write-Host "Report refreshed for handback of 7ed2a854-01c6-4bd3-8381-0ee20e055f16.md"
